$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.883.44'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '3.134.34'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.32'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.08'
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.128.67'
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("E9").Value = '  -0.22%  '
$ws.Range("E10").Value = '  -0.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.91'
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.34'
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("D15").Value = '3.653.39'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.32'
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("D18").Value = '63.750.18'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '3.133.32'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.84'
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.54'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.50'
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  +7.45%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  +10.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.47'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  +8.88%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.62'
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.20'
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = '  -6.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.27'
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = '  +0.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.34'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  +7.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '452.41'
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.292'
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").Value = '2.912.08'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("E46").Value = '  +12.20%  '
$ws.Range("E47").Value = '  -2.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.68'
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  +7.33%  '
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("E51").Value = '  +2.28%  '
